$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A2").Value = "Warringah Triathlon Club"
$ws.Range("B2").Value = 30
$ws.Range("C2").Value = 109
$ws.Range("D2").Value = 139
$ws.Range("E2").Value = 139
$ws.Range("F2").Value = 253
$ws.Range("A3").Value = "Balmoral Triathlon Club"
$ws.Range("B3").Value = 30
$ws.Range("C3").Value = 93
$ws.Range("D3").Value = 123
$ws.Range("E3").Value = 123
$ws.Range("F3").Value = 152
$ws.Range("A4").Value = "Panthers Triathlon Club"
$ws.Range("B4").Value = 30
$ws.Range("C4").Value = 79
$ws.Range("D4").Value = 109
$ws.Range("E4").Value = 109
$ws.Range("F4").Value = 135
$ws.Range("A5").Value = "Moore Performance Triathlon Club"
$ws.Range("B5").Value = 45
$ws.Range("C5").Value = 54
$ws.Range("D5").Value = 99
$ws.Range("E5").Value = 99
$ws.Range("F5").Value = 44
$ws.Range("A6").Value = "STG Triathlon Club"
$ws.Range("B6").Value = 45
$ws.Range("C6").Value = 47
$ws.Range("D6").Value = 92
$ws.Range("E6").Value = 92
$ws.Range("F6").Value = 49
$ws.Range("A7").Value = "Cronulla Triathlon Club"
$ws.Range("B7").Value = 15
$ws.Range("C7").Value = 65
$ws.Range("D7").Value = 80
$ws.Range("E7").Value = 80
$ws.Range("F7").Value = 207
$ws.Range("A8").Value = "Coogee Triathlon Club"
$ws.Range("B8").Value = 30
$ws.Range("C8").Value = 40
$ws.Range("D8").Value = 70
$ws.Range("E8").Value = 70
$ws.Range("F8").Value = 119
$ws.Range("A9").Value = "Hunters Hills Triathlon Club"
$ws.Range("B9").Value = 15
$ws.Range("C9").Value = 54
$ws.Range("D9").Value = 69
$ws.Range("E9").Value = 69
$ws.Range("F9").Value = 6
$ws.Range("A10").Value = "Pulse Performance"
$ws.Range("B10").Value = 30
$ws.Range("C10").Value = 26
$ws.Range("D10").Value = 56
$ws.Range("E10").Value = 56
$ws.Range("F10").Value = 31
$ws.Range("A11").Value = "Engadine Triathlon Club"
$ws.Range("B11").Value = 45
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 45
$ws.Range("E11").Value = 45
$ws.Range("F11").Value = 8
$ws.Range("A12").Value = "Concord Triathlon Club"
$ws.Range("B12").Value = 30
$ws.Range("C12").Value = 15
$ws.Range("D12").Value = 45
$ws.Range("E12").Value = 45
$ws.Range("F12").Value = 53
$ws.Range("A13").Value = "BRAT Triathlon Club"
$ws.Range("B13").Value = 30
$ws.Range("C13").Value = 9
$ws.Range("D13").Value = 39
$ws.Range("E13").Value = 39
$ws.Range("F13").Value = 161
$ws.Range("A14").Value = "Macarthur Triathlon Club"
$ws.Range("B14").Value = 15
$ws.Range("C14").Value = 19
$ws.Range("D14").Value = 34
$ws.Range("E14").Value = 34
$ws.Range("F14").Value = 69
$ws.Range("A15").Value = "Australian Chinese Dragon"
$ws.Range("B15").Value = 30
$ws.Range("C15").Value = 0
$ws.Range("D15").Value = 30
$ws.Range("E15").Value = 30
$ws.Range("F15").Value = 3
$ws.Range("A16").Value = "FilOz Triathlon Club"
$ws.Range("B16").Value = 30
$ws.Range("C16").Value = 0
$ws.Range("D16").Value = 30
$ws.Range("E16").Value = 30
$ws.Range("F16").Value = 10
$ws.Range("A17").Value = "South West Sydney Triathlon Club"
$ws.Range("B17").Value = 30
$ws.Range("C17").Value = 0
$ws.Range("D17").Value = 30
$ws.Range("E17").Value = 30
$ws.Range("F17").Value = 2
$ws.Range("A18").Value = "Northern Suburbs Triathlon Club"
$ws.Range("B18").Value = 15
$ws.Range("C18").Value = 9
$ws.Range("D18").Value = 24
$ws.Range("E18").Value = 24
$ws.Range("F18").Value = 67
$ws.Range("A19").Value = "Brighton Baths Athletic Club"
$ws.Range("B19").Value = 15
$ws.Range("C19").Value = 9
$ws.Range("D19").Value = 24
$ws.Range("E19").Value = 24
$ws.Range("F19").Value = 42
$ws.Range("A20").Value = "Manly Vipers Triathlon Club"
$ws.Range("B20").Value = 15
$ws.Range("C20").Value = 9
$ws.Range("D20").Value = 24
$ws.Range("E20").Value = 24
$ws.Range("F20").Value = 33
$ws.Range("A21").Value = "Balance Triathlon Club"
$ws.Range("B21").Value = 0
$ws.Range("C21").Value = 0
$ws.Range("D21").Value = 0
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 127
$ws.Range("A22").Value = "Hills Red Army"
$ws.Range("B22").Value = 0
$ws.Range("C22").Value = 0
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 178
$ws.Range("A23").Value = "TriMob"
$ws.Range("B23").Value = 0
$ws.Range("C23").Value = 0
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 43
